# Auto-generated edit script: applies cached market-data value updates
# to the Leve profit tables across all 8 class sheets (scheduled runner update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 379.66666
$ws.Range("I12").Value = 379.66666
$ws.Range("K12").Value = 379.66666
$ws.Range("M12").Value = -209.66666
$ws.Range("H43").Value = 3080.25
$ws.Range("I43").Value = 2948.6667
$ws.Range("J43").Value = 3475
$ws.Range("K43").Value = 2948.6667
$ws.Range("L43").Value = 3475
$ws.Range("M43").Value = -2879.6667
$ws.Range("N43").Value = -3613
$ws.Range("H76").Value = 4983.8965
$ws.Range("J76").Value = 5614.0557
$ws.Range("L76").Value = 5614.0557
$ws.Range("N76").Value = -6244.0557
$ws.Range("H79").Value = 4983.8965
$ws.Range("J79").Value = 5614.0557
$ws.Range("L79").Value = 5614.0557
$ws.Range("M79").Value = -2860.7273
$ws.Range("N79").Value = -7798.0557
$ws.Range("H80").Value = 793.2353000000001
$ws.Range("I80").Value = 1041.4445
$ws.Range("K80").Value = 3124.3335
$ws.Range("M80").Value = -2126.3335
$ws.Range("H83").Value = 793.2353000000001
$ws.Range("I83").Value = 1041.4445
$ws.Range("K83").Value = 9373.0005
$ws.Range("M83").Value = -4381.0005
$ws.Range("H86").Value = 6164.143
$ws.Range("I86").Value = 3537.25
$ws.Range("J86").Value = 9666.666999999999
$ws.Range("K86").Value = 3537.25
$ws.Range("L86").Value = 9666.666999999999
$ws.Range("M86").Value = -2414.25
$ws.Range("N86").Value = -11912.667
$ws.Range("H88").Value = 2894.6667
$ws.Range("I88").Value = 3586.25
$ws.Range("K88").Value = 3586.25
$ws.Range("M88").Value = -3180.25
$ws.Range("H89").Value = 6164.143
$ws.Range("I89").Value = 3537.25
$ws.Range("J89").Value = 9666.666999999999
$ws.Range("K89").Value = 17686.25
$ws.Range("L89").Value = 48333.335
$ws.Range("M89").Value = -12070.25
$ws.Range("N89").Value = -59565.335
$ws.Range("H91").Value = 2894.6667
$ws.Range("I91").Value = 3586.25
$ws.Range("K91").Value = 3586.25
$ws.Range("M91").Value = -2182.25
$ws.Range("H100").Value = 1157.2858
$ws.Range("I100").Value = 933.5
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 933.5
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -392.5
$ws.Range("N100").Value = -3582
$ws.Range("H103").Value = 1646.7894
$ws.Range("I103").Value = 1605.9333
$ws.Range("J103").Value = 1800
$ws.Range("K103").Value = 4817.7999
$ws.Range("L103").Value = 5400
$ws.Range("M103").Value = -4231.7999
$ws.Range("N103").Value = -6572
$ws.Range("H106").Value = 12110.5
$ws.Range("I106").Value = 7376.8
$ws.Range("K106").Value = 7376.8
$ws.Range("M106").Value = -6745.8
$ws.Range("H107").Value = 3846.9167
$ws.Range("I107").Value = 2221.1428
$ws.Range("K107").Value = 2221.1428
$ws.Range("M107").Value = -301.1428000000001
$ws.Range("H113").Value = 5242.3335
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = ""
$ws.Range("H116").Value = 22588.5
$ws.Range("I116").Value = 10660.111
$ws.Range("J116").Value = 44059.6
$ws.Range("K116").Value = 10660.111
$ws.Range("L116").Value = 44059.6
$ws.Range("M116").Value = -7218.111000000001
$ws.Range("N116").Value = -50943.6
$ws.Range("H125").Value = 1997
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1997
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 17973
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = -22893
$ws.Range("H127").Value = 1102575.6
$ws.Range("I127").Value = 1574039.4
$ws.Range("J127").Value = 2493.3333
$ws.Range("K127").Value = 4722118.199999999
$ws.Range("L127").Value = 7479.999899999999
$ws.Range("M127").Value = -4717158.199999999
$ws.Range("N127").Value = -17399.9999
$ws.Range("H131").Value = 1251235.9
$ws.Range("I131").Value = 1514.6666
$ws.Range("J131").Value = 5000399.5
$ws.Range("K131").Value = 4543.9998
$ws.Range("L131").Value = 15001198.5
$ws.Range("M131").Value = 496.0002000000004
$ws.Range("N131").Value = -15011278.5
$ws.Range("H132").Value = 1654.8235
$ws.Range("I132").Value = 1622.8667
$ws.Range("K132").Value = 4868.6001
$ws.Range("M132").Value = -2338.6001
$ws.Range("H135").Value = 3849.0322
$ws.Range("I135").Value = 2281.5386
$ws.Range("K135").Value = 20533.8474
$ws.Range("M135").Value = -17998.8474
$ws.Range("H137").Value = 927.23254
$ws.Range("I137").Value = 852.9459000000001
$ws.Range("J137").Value = 1385.3334
$ws.Range("K137").Value = 2558.8377
$ws.Range("L137").Value = 4156.0002
$ws.Range("M137").Value = -8.837700000000041
$ws.Range("N137").Value = -9256.0002
$ws.Range("H141").Value = 4145.4814
$ws.Range("I141").Value = 1773.15
$ws.Range("J141").Value = 10923.571
$ws.Range("K141").Value = 5319.450000000001
$ws.Range("L141").Value = 32770.713
$ws.Range("M141").Value = -139.4500000000007
$ws.Range("N141").Value = -43130.713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 488.14285
$ws.Range("I2").Value = 488.14285
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 488.14285
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -375.14285
$ws.Range("N2").Value = ""
$ws.Range("H14").Value = 1265.4
$ws.Range("I14").Value = 120
$ws.Range("J14").Value = 1551.75
$ws.Range("K14").Value = 120
$ws.Range("L14").Value = 1551.75
$ws.Range("M14").Value = 55
$ws.Range("N14").Value = -1901.75
$ws.Range("H16").Value = 499
$ws.Range("I16").Value = 496.33334
$ws.Range("J16").Value = 507
$ws.Range("K16").Value = 496.33334
$ws.Range("L16").Value = 507
$ws.Range("M16").Value = -209.33334
$ws.Range("N16").Value = -1081
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""
$ws.Range("H32").Value = 817.5139
$ws.Range("I32").Value = 627.0161000000001
$ws.Range("K32").Value = 627.0161000000001
$ws.Range("M32").Value = -340.0161000000001
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = ""
$ws.Range("H45").Value = 3169.8
$ws.Range("I45").Value = 4100
$ws.Range("J45").Value = 1774.5
$ws.Range("K45").Value = 4100
$ws.Range("L45").Value = 1774.5
$ws.Range("M45").Value = -3723
$ws.Range("N45").Value = -2528.5
$ws.Range("H97").Value = 537.67645
$ws.Range("I97").Value = 559.4375
$ws.Range("J97").Value = 189.5
$ws.Range("K97").Value = 559.4375
$ws.Range("L97").Value = 189.5
$ws.Range("M97").Value = -63.4375
$ws.Range("N97").Value = -1181.5
$ws.Range("H111").Value = 65000
$ws.Range("J111").Value = 65000
$ws.Range("L111").Value = 65000
$ws.Range("N111").Value = -73180
$ws.Range("H116").Value = 488.14285
$ws.Range("I116").Value = 488.14285
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 488.14285
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1805.85715
$ws.Range("N116").Value = ""
$ws.Range("H131").Value = 90712
$ws.Range("J131").Value = 90712
$ws.Range("L131").Value = 90712
$ws.Range("N131").Value = -100792
$ws.Range("H132").Value = 5822.8823
$ws.Range("I132").Value = 4377.2856
$ws.Range("K132").Value = 13131.8568
$ws.Range("M132").Value = -10601.8568
$ws.Range("H138").Value = 99996.336
$ws.Range("J138").Value = 99996.336
$ws.Range("L138").Value = 99996.336
$ws.Range("N138").Value = -110276.336
$ws.Range("H139").Value = 86874.86
$ws.Range("J139").Value = 86874.86
$ws.Range("L139").Value = 86874.86
$ws.Range("N139").Value = -97154.86

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 488.14285
$ws.Range("I3").Value = 488.14285
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 488.14285
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -374.14285
$ws.Range("N3").Value = ""
$ws.Range("H20").Value = 7201.4375
$ws.Range("I20").Value = 5025.385
$ws.Range("K20").Value = 5025.385
$ws.Range("M20").Value = -4778.385
$ws.Range("H35").Value = 35333
$ws.Range("J35").Value = 35333
$ws.Range("L35").Value = 35333
$ws.Range("N35").Value = -35953
$ws.Range("H64").Value = 9206.916999999999
$ws.Range("I64").Value = 14686.286
$ws.Range("K64").Value = 14686.286
$ws.Range("M64").Value = -14461.286
$ws.Range("H67").Value = 9206.916999999999
$ws.Range("I67").Value = 14686.286
$ws.Range("K67").Value = 14686.286
$ws.Range("M67").Value = -13906.286
$ws.Range("H94").Value = 3800.4167
$ws.Range("I94").Value = 3760.5
$ws.Range("K94").Value = 3760.5
$ws.Range("M94").Value = -3309.5
$ws.Range("H105").Value = 4614
$ws.Range("I105").Value = 3644.8333
$ws.Range("J105").Value = 5583.1665
$ws.Range("K105").Value = 3644.8333
$ws.Range("L105").Value = 5583.1665
$ws.Range("M105").Value = -1897.8333
$ws.Range("N105").Value = -9077.166499999999
$ws.Range("H107").Value = 1336.2222
$ws.Range("I107").Value = 1380.6666
$ws.Range("J107").Value = 1280.6666
$ws.Range("K107").Value = 1380.6666
$ws.Range("L107").Value = 1280.6666
$ws.Range("M107").Value = 539.3334
$ws.Range("N107").Value = -5120.6666
$ws.Range("H123").Value = 82500
$ws.Range("J123").Value = 82500
$ws.Range("L123").Value = 82500
$ws.Range("N123").Value = -92300
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""
$ws.Range("H134").Value = 2421.068
$ws.Range("I134").Value = 1720.2927
$ws.Range("K134").Value = 5160.8781
$ws.Range("M134").Value = -2625.8781

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 625325
$ws.Range("I4").Value = 2500000
$ws.Range("J4").Value = 433.33334
$ws.Range("K4").Value = 2500000
$ws.Range("L4").Value = 433.33334
$ws.Range("M4").Value = -2499888
$ws.Range("N4").Value = -657.33334
$ws.Range("H14").Value = 17669
$ws.Range("I14").Value = 25754
$ws.Range("J14").Value = 1499
$ws.Range("K14").Value = 25754
$ws.Range("L14").Value = 1499
$ws.Range("M14").Value = -25584
$ws.Range("N14").Value = -1839
$ws.Range("H16").Value = 69079.07000000001
$ws.Range("I16").Value = 2277.4443
$ws.Range("J16").Value = 169281.5
$ws.Range("K16").Value = 2277.4443
$ws.Range("L16").Value = 169281.5
$ws.Range("M16").Value = -1990.4443
$ws.Range("N16").Value = -169855.5
$ws.Range("H21").Value = 600
$ws.Range("I21").Value = 600
$ws.Range("K21").Value = 600
$ws.Range("M21").Value = -365
$ws.Range("H25").Value = 843.8570999999999
$ws.Range("I25").Value = 843.8570999999999
$ws.Range("K25").Value = 843.8570999999999
$ws.Range("M25").Value = -669.8570999999999
$ws.Range("H31").Value = 275485.72
$ws.Range("I31").Value = 831709
$ws.Range("J31").Value = 4890.6216
$ws.Range("K31").Value = 831709
$ws.Range("L31").Value = 4890.6216
$ws.Range("M31").Value = -831414
$ws.Range("N31").Value = -5480.6216
$ws.Range("H34").Value = 275485.72
$ws.Range("I34").Value = 831709
$ws.Range("J34").Value = 4890.6216
$ws.Range("K34").Value = 831709
$ws.Range("L34").Value = 4890.6216
$ws.Range("M34").Value = -831507
$ws.Range("N34").Value = -5294.6216
$ws.Range("H58").Value = 4077.261
$ws.Range("I58").Value = 7833.3335
$ws.Range("J58").Value = 3513.85
$ws.Range("K58").Value = 7833.3335
$ws.Range("L58").Value = 3513.85
$ws.Range("M58").Value = -7630.3335
$ws.Range("N58").Value = -3919.85
$ws.Range("H113").Value = 69079.07000000001
$ws.Range("I113").Value = 2277.4443
$ws.Range("J113").Value = 169281.5
$ws.Range("K113").Value = 2277.4443
$ws.Range("L113").Value = 169281.5
$ws.Range("M113").Value = -107.4443000000001
$ws.Range("N113").Value = -173621.5
$ws.Range("H122").Value = 12944.632
$ws.Range("I122").Value = 2272.923
$ws.Range("K122").Value = 6818.768999999999
$ws.Range("M122").Value = -4368.768999999999
$ws.Range("H132").Value = 3176.7568
$ws.Range("I132").Value = 3037.84
$ws.Range("J132").Value = 3466.1667
$ws.Range("K132").Value = 9113.52
$ws.Range("L132").Value = 10398.5001
$ws.Range("M132").Value = -6583.52
$ws.Range("N132").Value = -15458.5001
$ws.Range("H134").Value = 3963.5305
$ws.Range("I134").Value = 2219.5789
$ws.Range("J134").Value = 5068.033
$ws.Range("K134").Value = 6658.736699999999
$ws.Range("L134").Value = 15204.099
$ws.Range("M134").Value = -4123.736699999999
$ws.Range("N134").Value = -20274.099
$ws.Range("H136").Value = 4077.261
$ws.Range("I136").Value = 7833.3335
$ws.Range("J136").Value = 3513.85
$ws.Range("K136").Value = 23500.0005
$ws.Range("L136").Value = 10541.55
$ws.Range("M136").Value = -20950.0005
$ws.Range("N136").Value = -15641.55

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I11").Value = 142.25
$ws.Range("K11").Value = 426.75
$ws.Range("M11").Value = -286.75
$ws.Range("H34").Value = 781.55554
$ws.Range("I34").Value = 754.25
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 2262.75
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -2178.75
$ws.Range("N34").Value = -3168
$ws.Range("H68").Value = 5130.78
$ws.Range("J68").Value = 5820.3022
$ws.Range("L68").Value = 17460.9066
$ws.Range("N68").Value = -19082.9066
$ws.Range("H71").Value = 5130.78
$ws.Range("J71").Value = 5820.3022
$ws.Range("L71").Value = 52382.7198
$ws.Range("N71").Value = -60494.7198
$ws.Range("H131").Value = 3548.318
$ws.Range("I131").Value = 1312.5385
$ws.Range("J131").Value = 6777.778
$ws.Range("K131").Value = 3937.6155
$ws.Range("L131").Value = 20333.334
$ws.Range("M131").Value = 1102.3845
$ws.Range("N131").Value = -30413.334
$ws.Range("H133").Value = 11347.467
$ws.Range("I133").Value = 3844
$ws.Range("K133").Value = 11532
$ws.Range("M133").Value = -6472
$ws.Range("H134").Value = 7919.4375
$ws.Range("I134").Value = 2745.6667
$ws.Range("K134").Value = 8237.000100000001
$ws.Range("M134").Value = -3167.000100000001
$ws.Range("H136").Value = 8357.546
$ws.Range("I136").Value = 4562
$ws.Range("J136").Value = 14999.75
$ws.Range("K136").Value = 13686
$ws.Range("L136").Value = 44999.25
$ws.Range("M136").Value = -8586
$ws.Range("N136").Value = -55199.25
$ws.Range("H138").Value = 18716.084
$ws.Range("I138").Value = 40478.75
$ws.Range("K138").Value = 121436.25
$ws.Range("M138").Value = -116296.25
$ws.Range("H139").Value = 6985.355
$ws.Range("I139").Value = 3347.4
$ws.Range("K139").Value = 10042.2
$ws.Range("M139").Value = -4902.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5050054.5
$ws.Range("I2").Value = 5941218.5
$ws.Range("K2").Value = 5941218.5
$ws.Range("M2").Value = -5941105.5
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("K5").Value = 500
$ws.Range("M5").Value = -388
$ws.Range("H21").Value = 625283.5
$ws.Range("I21").Value = 16000
$ws.Range("K21").Value = 16000
$ws.Range("M21").Value = -15827
$ws.Range("H30").Value = 625283.5
$ws.Range("I30").Value = 16000
$ws.Range("K30").Value = 16000
$ws.Range("M30").Value = -15895
$ws.Range("H97").Value = 1342.8572
$ws.Range("J97").Value = 1732.6666
$ws.Range("L97").Value = 1732.6666
$ws.Range("N97").Value = -2724.6666
$ws.Range("H102").Value = 1597.1923
$ws.Range("I102").Value = 1631.6957
$ws.Range("K102").Value = 1631.6957
$ws.Range("M102").Value = -9.695699999999988
$ws.Range("H126").Value = 2837.2222
$ws.Range("I126").Value = 2754.375
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 8263.125
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -5793.125
$ws.Range("N126").Value = -15440
$ws.Range("H132").Value = 17921.4
$ws.Range("I132").Value = 14136.695
$ws.Range("J132").Value = 30356.857
$ws.Range("K132").Value = 42410.085
$ws.Range("L132").Value = 91070.571
$ws.Range("M132").Value = -39880.085
$ws.Range("N132").Value = -96130.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3792.5334
$ws.Range("I7").Value = 3225
$ws.Range("K7").Value = 3225
$ws.Range("M7").Value = -3113
$ws.Range("H16").Value = 579.7273
$ws.Range("I16").Value = 598.9
$ws.Range("K16").Value = 598.9
$ws.Range("M16").Value = -428.9
$ws.Range("H22").Value = 2168.5745
$ws.Range("I22").Value = 1827.8
$ws.Range("J22").Value = 2555.818
$ws.Range("K22").Value = 1827.8
$ws.Range("L22").Value = 2555.818
$ws.Range("M22").Value = -1532.8
$ws.Range("N22").Value = -3145.818
$ws.Range("H23").Value = 795.5
$ws.Range("I23").Value = 795.5
$ws.Range("K23").Value = 795.5
$ws.Range("M23").Value = -565.5
$ws.Range("H25").Value = 700
$ws.Range("J25").Value = 700
$ws.Range("L25").Value = 700
$ws.Range("N25").Value = -1160
$ws.Range("H27").Value = 2168.5745
$ws.Range("I27").Value = 1827.8
$ws.Range("J27").Value = 2555.818
$ws.Range("K27").Value = 1827.8
$ws.Range("L27").Value = 2555.818
$ws.Range("M27").Value = -1720.8
$ws.Range("N27").Value = -2769.818
$ws.Range("H40").Value = 4154.154
$ws.Range("I40").Value = 3917.4167
$ws.Range("J40").Value = 6995
$ws.Range("K40").Value = 3917.4167
$ws.Range("L40").Value = 6995
$ws.Range("M40").Value = -3781.4167
$ws.Range("N40").Value = -7267
$ws.Range("H46").Value = 4388
$ws.Range("I46").Value = 3995
$ws.Range("J46").Value = 4411.1177
$ws.Range("K46").Value = 3995
$ws.Range("L46").Value = 4411.1177
$ws.Range("M46").Value = -3807
$ws.Range("N46").Value = -4787.1177
$ws.Range("H55").Value = 1257.9678
$ws.Range("I55").Value = 160.11111
$ws.Range("J55").Value = 1707.091
$ws.Range("K55").Value = 160.11111
$ws.Range("L55").Value = 1707.091
$ws.Range("M55").Value = 12.88889
$ws.Range("N55").Value = -2053.091
$ws.Range("H61").Value = 10383.4
$ws.Range("I61").Value = 10987.833
$ws.Range("J61").Value = 7965.6665
$ws.Range("K61").Value = 10987.833
$ws.Range("L61").Value = 7965.6665
$ws.Range("M61").Value = -10785.833
$ws.Range("N61").Value = -8369.666499999999
$ws.Range("H113").Value = 10383.4
$ws.Range("I113").Value = 10987.833
$ws.Range("J113").Value = 7965.6665
$ws.Range("K113").Value = 10987.833
$ws.Range("L113").Value = 7965.6665
$ws.Range("M113").Value = -8817.833000000001
$ws.Range("N113").Value = -12305.6665
$ws.Range("H117").Value = 80392
$ws.Range("J117").Value = 80392
$ws.Range("L117").Value = 80392
$ws.Range("N117").Value = -89570
$ws.Range("H122").Value = 4750
$ws.Range("J122").Value = 8000
$ws.Range("L122").Value = 24000
$ws.Range("N122").Value = -28900
$ws.Range("H126").Value = 3792.5334
$ws.Range("I126").Value = 3225
$ws.Range("K126").Value = 9675
$ws.Range("M126").Value = -7205
$ws.Range("H132").Value = 2635320.5
$ws.Range("I132").Value = 5558703.5
$ws.Range("J132").Value = 4275.9
$ws.Range("K132").Value = 16676110.5
$ws.Range("L132").Value = 12827.7
$ws.Range("M132").Value = -16673580.5
$ws.Range("N132").Value = -17887.7
$ws.Range("H134").Value = 83999.664
$ws.Range("J134").Value = 83999.664
$ws.Range("L134").Value = 83999.664
$ws.Range("N134").Value = -94139.664
$ws.Range("H136").Value = 3230396.5
$ws.Range("I136").Value = 4547833.5
$ws.Range("J136").Value = 9995.444
$ws.Range("K136").Value = 13643500.5
$ws.Range("L136").Value = 29986.332
$ws.Range("M136").Value = -13640950.5
$ws.Range("N136").Value = -35086.33199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 99000
$ws.Range("J16").Value = 99000
$ws.Range("L16").Value = 99000
$ws.Range("N16").Value = -99584
$ws.Range("H23").Value = 1582.3334
$ws.Range("I23").Value = 209.16667
$ws.Range("J23").Value = 4328.6665
$ws.Range("K23").Value = 209.16667
$ws.Range("L23").Value = 4328.6665
$ws.Range("M23").Value = 19.83332999999999
$ws.Range("N23").Value = -4786.6665
$ws.Range("H62").Value = 11720.615
$ws.Range("I62").Value = 4707.222
$ws.Range("K62").Value = 4707.222
$ws.Range("M62").Value = -4083.222
$ws.Range("H65").Value = 11720.615
$ws.Range("I65").Value = 4707.222
$ws.Range("K65").Value = 23536.11
$ws.Range("M65").Value = -20416.11
$ws.Range("H81").Value = 4363.636
$ws.Range("I81").Value = 3800
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 7600
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -6539
$ws.Range("N81").Value = -22122
$ws.Range("H84").Value = 4363.636
$ws.Range("I84").Value = 3800
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 38000
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -32696
$ws.Range("N84").Value = -110608
$ws.Range("H122").Value = 41035.516
$ws.Range("I122").Value = 2101.5557
$ws.Range("K122").Value = 6304.6671
$ws.Range("M122").Value = -3854.6671
$ws.Range("H124").Value = 50139.332
$ws.Range("J124").Value = 50139.332
$ws.Range("L124").Value = 50139.332
$ws.Range("N124").Value = -59959.332
$ws.Range("H132").Value = 3766.9285
$ws.Range("I132").Value = 3519.75
$ws.Range("K132").Value = 10559.25
$ws.Range("M132").Value = -8029.25
$ws.Range("H135").Value = 52898
$ws.Range("J135").Value = 52898
$ws.Range("L135").Value = 52898
$ws.Range("N135").Value = -63038
$ws.Range("H136").Value = 41669010
$ws.Range("I136").Value = 50002060
$ws.Range("J136").Value = 3747.25
$ws.Range("K136").Value = 150006180
$ws.Range("L136").Value = 11241.75
$ws.Range("M136").Value = -150003630
$ws.Range("N136").Value = -16341.75
